$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (issue number and week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Column E width (matches other narrow columns; closest attainable to target) ---
$ws.Columns.Item(5).ColumnWidth = 5.43

# --- Cells changing number format style (copy format from a same-style cell, then set value) ---
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("I14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2

# --- Cells changing to text placeholders (force text format, set value, restore General style) ---
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -88
$ws.Range("F15").Value = 9
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 350
$ws.Range("I15").Value = 23
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 155.555555555556
$ws.Range("L15").Value = 27.777777777777
$ws.Range("M15").Value = 130
$ws.Range("N15").Value = -4.166666666666
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = 8.333333333333
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 40
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 131
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = -14.935064935064
$ws.Range("L16").Value = -23.837209302325
$ws.Range("M16").Value = -9.027777777777
$ws.Range("N16").Value = -79.658385093167
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 17.647058823529
$ws.Range("I17").Value = 232
$ws.Range("J17").Value = 207
$ws.Range("K17").Value = 12.077294685990
$ws.Range("L17").Value = -10.424710424710
$ws.Range("M17").Value = 45
$ws.Range("N17").Value = -10.769230769230
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 5.882352941176
$ws.Range("I18").Value = 90
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = -5.263157894736
$ws.Range("L18").Value = -28
$ws.Range("M18").Value = -26.229508196721
$ws.Range("N18").Value = -83.754512635379
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 91.666666666666
$ws.Range("F19").Value = 99
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = 62.295081967213
$ws.Range("I19").Value = 371
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = 31.095406360424
$ws.Range("L19").Value = 17.777777777777
$ws.Range("M19").Value = 130.434782608696
$ws.Range("N19").Value = 53.941908713692
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -15.384615384615
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -3.030303030303
$ws.Range("I20").Value = 145
$ws.Range("J20").Value = 147
$ws.Range("K20").Value = -1.360544217687
$ws.Range("L20").Value = -32.242990654205
$ws.Range("M20").Value = 123.076923076923
$ws.Range("N20").Value = -76.837060702875
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = 22.641509433962
$ws.Range("F21").Value = 246
$ws.Range("G21").Value = 204
$ws.Range("H21").Value = 20.588235294117
$ws.Range("I21").Value = 995
$ws.Range("J21").Value = 896
$ws.Range("K21").Value = 11.049107142857
$ws.Range("L21").Value = -10.036166365280
$ws.Range("M21").Value = 49.849397590361
$ws.Range("N21").Value = -58.087615838247
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -55.555555555555
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -8
$ws.Range("I23").Value = 77
$ws.Range("J23").Value = 85
$ws.Range("K23").Value = -9.411764705882
$ws.Range("L23").Value = -14.444444444444
$ws.Range("M23").Value = 13.235294117647
$ws.Range("C24").Value = 56
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = 51.351351351351
$ws.Range("F24").Value = 171
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = 56.880733944954
$ws.Range("I24").Value = 711
$ws.Range("J24").Value = 549
$ws.Range("K24").Value = 29.508196721311
$ws.Range("L24").Value = 12.678288431061
$ws.Range("M24").Value = 61.958997722095
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 54.285714285714
$ws.Range("I25").Value = 226
$ws.Range("J25").Value = 192
$ws.Range("K25").Value = 17.708333333333
$ws.Range("L25").Value = -17.518248175182
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = 4.545454545454
$ws.Range("F26").Value = 93
$ws.Range("G26").Value = 80
$ws.Range("H26").Value = 16.25
$ws.Range("I26").Value = 332
$ws.Range("J26").Value = 369
$ws.Range("K26").Value = -10.027100271002
$ws.Range("L26").Value = -0.895522388059
$ws.Range("M26").Value = -29.059829059829
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 25
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 47.058823529411
$ws.Range("L27").Value = 4.166666666666
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -42.857142857142
$ws.Range("F29").Value = 2
$ws.Range("N29").Value = -80
$ws.Range("F30").Value = 2
$ws.Range("N30").Value = -78.571428571428
